# Created functions to get season record
# Add the team's season Wins/Losses/Ties as three new trailing columns
# (AD, AE, AF) on the player-stats sheet, matching the header styling
# already used by the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold,
# centered, bordered) onto the three new header cells, then set text ---
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-56): season record is a team-level stat, so every
# player row gets the same Wins/Losses/Ties values ---
$ws.Range("AD2:AD56").Value = 93
$ws.Range("AE2:AE56").Value = 69
$ws.Range("AF2:AF56").Value = 0
